# Remove the "oldpeak" (E) and "thall" (F) columns from the normalized
# data sheet, shrinking the used range from A1:F304 down to A1:D304.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1:F1").EntireColumn.Delete()
